$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width (~20.7 characters)
$ws.Columns.Item(1).ColumnWidth = 19.833333333333336

# Center-align A1:B2
$ws.Range("A1:B2").HorizontalAlignment = -4108

# C1:D2 bold + centered, then merge with D1:D2
$ws.Range("C1:D2").HorizontalAlignment = -4108
$ws.Range("C1:D2").Font.Bold = $true

$ws.Range("C2").Value = "1. siječnja 2023."

$ws.Range("C1:D1").Merge()
$ws.Range("C2:D2").Merge()

$ws.Range("D4").Select()
